# Weekly data refresh: insert 3 new observation rows at the top of the
# "Ajo" data block (original rows 44-136 shift down to 47-139), and
# populate the 3 new rows (44-46) with the latest weekly prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 44; existing rows 44-136 become 47-139.
$ws.Rows("44:46").Insert()

# Shared values for the new rows (same market/category across the sheet).
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112003
$categoria = "Ajo"
$clasif    = "Hortaliza"

# New row 44: Rosado / 1a nueva(o)
$ws.Cells.Item(44, 1).Value  = $mercadoId
$ws.Cells.Item(44, 2).Value  = $mercado
$ws.Cells.Item(44, 3).Value  = $region
$ws.Cells.Item(44, 4).Value  = 44519
$ws.Cells.Item(44, 5).Value  = $codreg
$ws.Cells.Item(44, 6).Value  = $catId
$ws.Cells.Item(44, 7).Value  = $categoria
$ws.Cells.Item(44, 8).Value  = "Rosado"
$ws.Cells.Item(44, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(44, 10).Value = 680
$ws.Cells.Item(44, 11).Value = 2200
$ws.Cells.Item(44, 12).Value = 2600
$ws.Cells.Item(44, 13).Value = 2400
$ws.Cells.Item(44, 14).Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(44, 15).Value = "Provincia de Talagante"
$ws.Cells.Item(44, 16).Value = 120
$ws.Cells.Item(44, 17).Value = 20
$ws.Cells.Item(44, 18).Value = $clasif

# New row 45: Rosado / 2a nueva(o)
$ws.Cells.Item(45, 1).Value  = $mercadoId
$ws.Cells.Item(45, 2).Value  = $mercado
$ws.Cells.Item(45, 3).Value  = $region
$ws.Cells.Item(45, 4).Value  = 44519
$ws.Cells.Item(45, 5).Value  = $codreg
$ws.Cells.Item(45, 6).Value  = $catId
$ws.Cells.Item(45, 7).Value  = $categoria
$ws.Cells.Item(45, 8).Value  = "Rosado"
$ws.Cells.Item(45, 9).Value  = "2a nueva(o)"
$ws.Cells.Item(45, 10).Value = 320
$ws.Cells.Item(45, 11).Value = 2000
$ws.Cells.Item(45, 12).Value = 2000
$ws.Cells.Item(45, 13).Value = 2000
$ws.Cells.Item(45, 14).Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(45, 15).Value = "Provincia de Talagante"
$ws.Cells.Item(45, 16).Value = 100
$ws.Cells.Item(45, 17).Value = 20
$ws.Cells.Item(45, 18).Value = $clasif

# New row 46: Rosado / Extra nueva (o)
$ws.Cells.Item(46, 1).Value  = $mercadoId
$ws.Cells.Item(46, 2).Value  = $mercado
$ws.Cells.Item(46, 3).Value  = $region
$ws.Cells.Item(46, 4).Value  = 44519
$ws.Cells.Item(46, 5).Value  = $codreg
$ws.Cells.Item(46, 6).Value  = $catId
$ws.Cells.Item(46, 7).Value  = $categoria
$ws.Cells.Item(46, 8).Value  = "Rosado"
$ws.Cells.Item(46, 9).Value  = "Extra nueva (o)"
$ws.Cells.Item(46, 10).Value = 500
$ws.Cells.Item(46, 11).Value = 2800
$ws.Cells.Item(46, 12).Value = 3000
$ws.Cells.Item(46, 13).Value = 2900
$ws.Cells.Item(46, 14).Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(46, 15).Value = "Provincia de Talagante"
$ws.Cells.Item(46, 16).Value = 145
$ws.Cells.Item(46, 17).Value = 20
$ws.Cells.Item(46, 18).Value = $clasif

# Dimension now spans through row 139.
Write-Output ("Rows: " + $ws.UsedRange.Rows.Count)
